$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 1.6
$ws.Cells.Item(3, 8).Value = 3.45
$ws.Cells.Item(3, 9).Value = 6.2
$ws.Cells.Item(3, 11).Value = 5.7
$ws.Cells.Item(3, 13).Value = 2.5
$ws.Cells.Item(3, 14).Value = 2.37
$ws.Cells.Item(3, 18).Value = 2.37
$ws.Cells.Item(3, 19).Value = 1.52
$ws.Cells.Item(3, 20).Value = 4.75
$ws.Cells.Item(3, 21).Value = 6
$ws.Cells.Item(3, 22).Value = 9
$ws.Cells.Item(3, 23).Value = 10.75
$ws.Cells.Item(3, 24).Value = 16
$ws.Cells.Item(3, 25).Value = 45
$ws.Cells.Item(3, 26).Value = 5.7
$ws.Cells.Item(3, 27).Value = 7.1
$ws.Cells.Item(3, 28).Value = 25
$ws.Cells.Item(3, 29).Value = 175
$ws.Cells.Item(3, 30).Value = 11.5
$ws.Cells.Item(3, 31).Value = 35
$ws.Cells.Item(3, 32).Value = 22
$ws.Cells.Item(3, 33).Value = 150
$ws.Cells.Item(3, 34).Value = 90
$ws.Cells.Item(3, 35).Value = 110

# Row 4
$ws.Cells.Item(4, 7).Value = 1.9
$ws.Cells.Item(4, 8).Value = 3.1
$ws.Cells.Item(4, 9).Value = 4.35
$ws.Cells.Item(4, 13).Value = 2.57
$ws.Cells.Item(4, 15).Value = 1.55
$ws.Cells.Item(4, 20).Value = 5.6
$ws.Cells.Item(4, 21).Value = 7.9
$ws.Cells.Item(4, 24).Value = 17.5
$ws.Cells.Item(4, 27).Value = 6.1
$ws.Cells.Item(4, 28).Value = 18
$ws.Cells.Item(4, 29).Value = 110
$ws.Cells.Item(4, 30).Value = 9.5
$ws.Cells.Item(4, 31).Value = 22
$ws.Cells.Item(4, 34).Value = 50
$ws.Cells.Item(4, 35).Value = 60

# Row 9
$ws.Cells.Item(9, 7).Value = 1.53
$ws.Cells.Item(9, 8).Value = 4.35
$ws.Cells.Item(9, 9).Value = 4.8
$ws.Cells.Item(9, 14).Value = 1.47
$ws.Cells.Item(9, 15).Value = 2.35
$ws.Cells.Item(9, 20).Value = 9.75
$ws.Cells.Item(9, 21).Value = 9
$ws.Cells.Item(9, 23).Value = 12
$ws.Cells.Item(9, 24).Value = 11
$ws.Cells.Item(9, 25).Value = 19.5
$ws.Cells.Item(9, 27).Value = 9
$ws.Cells.Item(9, 28).Value = 14.5
$ws.Cells.Item(9, 30).Value = 18.5
$ws.Cells.Item(9, 31).Value = 32
$ws.Cells.Item(9, 32).Value = 16
$ws.Cells.Item(9, 33).Value = 80
$ws.Cells.Item(9, 34).Value = 40
$ws.Cells.Item(9, 35).Value = 37

# Row 10
$ws.Cells.Item(10, 8).Value = 3.7
$ws.Cells.Item(10, 9).Value = 3.15
$ws.Cells.Item(10, 14).Value = 1.52
$ws.Cells.Item(10, 15).Value = 2.22
$ws.Cells.Item(10, 19).Value = 2.32
$ws.Cells.Item(10, 20).Value = 10.25
$ws.Cells.Item(10, 21).Value = 11.75
$ws.Cells.Item(10, 22).Value = 8.5
$ws.Cells.Item(10, 25).Value = 19.5
$ws.Cells.Item(10, 26).Value = 15.5
$ws.Cells.Item(10, 27).Value = 7.6
$ws.Cells.Item(10, 28).Value = 11.75
$ws.Cells.Item(10, 29).Value = 40
$ws.Cells.Item(10, 30).Value = 14
$ws.Cells.Item(10, 31).Value = 21
$ws.Cells.Item(10, 33).Value = 45
$ws.Cells.Item(10, 34).Value = 24
$ws.Cells.Item(10, 35).Value = 25

# Row 11
$ws.Cells.Item(11, 8).Value = 4.8
$ws.Cells.Item(11, 9).Value = 5.9
$ws.Cells.Item(11, 14).Value = 1.36
$ws.Cells.Item(11, 17).Value = 3.9
$ws.Cells.Item(11, 18).Value = 1.5
$ws.Cells.Item(11, 19).Value = 2.42
$ws.Cells.Item(11, 20).Value = 12.5
$ws.Cells.Item(11, 21).Value = 10
$ws.Cells.Item(11, 23).Value = 11.75
$ws.Cells.Item(11, 27).Value = 10.5
$ws.Cells.Item(11, 30).Value = 26
$ws.Cells.Item(11, 31).Value = 45
$ws.Cells.Item(11, 32).Value = 19

# Row 12
$ws.Cells.Item(12, 7).Value = 2.1
$ws.Cells.Item(12, 8).Value = 2.9
$ws.Cells.Item(12, 9).Value = 3.6
$ws.Cells.Item(12, 11).Value = 6.1
$ws.Cells.Item(12, 12).Value = 1.39
$ws.Cells.Item(12, 13).Value = 2.77
$ws.Cells.Item(12, 14).Value = 2.15
$ws.Cells.Item(12, 18).Value = 1.83
$ws.Cells.Item(12, 19).Value = 1.87
$ws.Cells.Item(12, 20).Value = 6.7
$ws.Cells.Item(12, 21).Value = 9.75
$ws.Cells.Item(12, 22).Value = 8.5
$ws.Cells.Item(12, 23).Value = 21
$ws.Cells.Item(12, 24).Value = 18
$ws.Cells.Item(12, 25).Value = 30
$ws.Cells.Item(12, 26).Value = 6.1
$ws.Cells.Item(12, 27).Value = 5.7
$ws.Cells.Item(12, 28).Value = 14.5
$ws.Cells.Item(12, 29).Value = 75
$ws.Cells.Item(12, 31).Value = 19.5
$ws.Cells.Item(12, 32).Value = 12.5
$ws.Cells.Item(12, 33).Value = 60
$ws.Cells.Item(12, 34).Value = 37
$ws.Cells.Item(12, 36).Value = 600

# Row 14
$ws.Cells.Item(14, 7).Value = 1.8
$ws.Cells.Item(14, 9).Value = 4.45
$ws.Cells.Item(14, 12).Value = 1.5
$ws.Cells.Item(14, 13).Value = 2.27
$ws.Cells.Item(14, 14).Value = 2.42
$ws.Cells.Item(14, 15).Value = 1.44
$ws.Cells.Item(14, 16).Value = 1.55
$ws.Cells.Item(14, 17).Value = 2.15
$ws.Cells.Item(14, 18).Value = 2.25
$ws.Cells.Item(14, 19).Value = 1.5
$ws.Cells.Item(14, 20).Value = 5
$ws.Cells.Item(14, 21).Value = 6.9
$ws.Cells.Item(14, 23).Value = 14
$ws.Cells.Item(14, 24).Value = 18.5
$ws.Cells.Item(14, 25).Value = 45
$ws.Cells.Item(14, 26).Value = 6.6
$ws.Cells.Item(14, 27).Value = 6.6
$ws.Cells.Item(14, 28).Value = 23
$ws.Cells.Item(14, 29).Value = 175
$ws.Cells.Item(14, 30).Value = 8.75
$ws.Cells.Item(14, 31).Value = 22
$ws.Cells.Item(14, 32).Value = 16.5
$ws.Cells.Item(14, 33).Value = 80
$ws.Cells.Item(14, 34).Value = 60
$ws.Cells.Item(14, 35).Value = 80

# Row 15
$ws.Cells.Item(15, 7).Value = 2.47
$ws.Cells.Item(15, 8).Value = 3
$ws.Cells.Item(15, 9).Value = 2.85
$ws.Cells.Item(15, 12).Value = 1.53
$ws.Cells.Item(15, 13).Value = 2.18
$ws.Cells.Item(15, 14).Value = 2.52
$ws.Cells.Item(15, 15).Value = 1.4
$ws.Cells.Item(15, 16).Value = 1.57
$ws.Cells.Item(15, 17).Value = 2.1
$ws.Cells.Item(15, 18).Value = 2.15
$ws.Cells.Item(15, 19).Value = 1.55
$ws.Cells.Item(15, 20).Value = 5.9
$ws.Cells.Item(15, 21).Value = 10.25
$ws.Cells.Item(15, 22).Value = 10.75
$ws.Cells.Item(15, 23).Value = 26
$ws.Cells.Item(15, 24).Value = 27
$ws.Cells.Item(15, 25).Value = 50
$ws.Cells.Item(15, 26).Value = 6.2
$ws.Cells.Item(15, 27).Value = 6.1
$ws.Cells.Item(15, 28).Value = 21
$ws.Cells.Item(15, 29).Value = 150
$ws.Cells.Item(15, 30).Value = 6.4
$ws.Cells.Item(15, 31).Value = 12.5
$ws.Cells.Item(15, 32).Value = 11.75
$ws.Cells.Item(15, 33).Value = 35
$ws.Cells.Item(15, 34).Value = 32
$ws.Cells.Item(15, 35).Value = 55

# Row 16
$ws.Cells.Item(16, 7).Value = 1.85
$ws.Cells.Item(16, 9).Value = 4
$ws.Cells.Item(16, 10).Value = 1.08
$ws.Cells.Item(16, 11).Value = 8
$ws.Cells.Item(16, 12).Value = 1.4
$ws.Cells.Item(16, 13).Value = 2.75
$ws.Cells.Item(16, 15).Value = 1.62
$ws.Cells.Item(16, 16).Value = 1.44
$ws.Cells.Item(16, 17).Value = 2.63
$ws.Cells.Item(16, 18).Value = 2
$ws.Cells.Item(16, 19).Value = 1.73
$ws.Cells.Item(16, 20).Value = 6
$ws.Cells.Item(16, 21).Value = 8
$ws.Cells.Item(16, 22).Value = 9
$ws.Cells.Item(16, 23).Value = 15
$ws.Cells.Item(16, 24).Value = 17
$ws.Cells.Item(16, 25).Value = 34
$ws.Cells.Item(16, 26).Value = 8
$ws.Cells.Item(16, 27).Value = 6.5
$ws.Cells.Item(16, 28).Value = 19
$ws.Cells.Item(16, 29).Value = 67
$ws.Cells.Item(16, 30).Value = 10
$ws.Cells.Item(16, 31).Value = 21
$ws.Cells.Item(16, 33).Value = 41
$ws.Cells.Item(16, 34).Value = 41
$ws.Cells.Item(16, 35).Value = 41
